$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-18 10:32:51"
$wsZh.Range("H3").Value = "2016-03-18 10:33:12"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-18 10:32:53"
$wsDe.Range("H3").Value = "2016-03-18 10:33:17"
